# Placement.xlsx edit script
# Summary of change (from commit "Add files via upload"):
#   On the "ProviderSearch" sheet, two new columns (CITY, ZIP_CODE) were
#   inserted right after the existing ADDRESS_LINE1 column (column Q),
#   pushing the validate/save/email/etc. columns two places to the right.
#   The demo/sample row (row 5) was filled in with a concrete city/zip
#   ("Boardman" / 97818) plus a couple of "Click" placeholders, and the
#   active selection/tab moved to the ProviderSearch sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProviderSearch")

# Insert 2 new blank columns before column R (pushes R:X -> T:Z).
$ws.Columns("R:S").Insert()

# New column headers in row 1.
$ws.Range("R1").Value = "CITY"
$ws.Range("S1").Value = "ZIP_CODE"

# Fill in the sample data on row 5 (the row that already has real data).
$ws.Range("Q5").Value = "autoText"
$ws.Range("R5").Value = "Boardman"
$ws.Range("S5").Value = 97818
$ws.Range("V5").Value = "Click"
$ws.Range("W5").Value = "Click"

# Make ProviderSearch the active sheet / tab, with the new columns selected.
$ws.Activate()
$ws.Range("S5").Select()

# The FolioChildLocations sheet view scrolled and the selection moved.
$ws2 = $wb.Worksheets.Item("FolioChildLocations")
$ws2.Range("R7").Select()

# Re-activate ProviderSearch so it stays the last-active (and therefore
# the tab persisted in the saved workbook view).
$ws.Activate()
